# Applies the "Quantum Entanglement" -> "Mathematics" essay rewrite.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find.Execute could not locate: $find"
    }
}

# --- Title / byline / contact block -----------------------------------
Replace-Text "Quantum Entanglement: A Tapestry of Interconnectedness" "Mathematics: A Journey Through Numbers and Patterns"
Replace-Text " Isabella Wilkinson" " Alexandro Alfonso Cintron"
Replace-Text "isabella" "aacintron@highlandschool"
Replace-Text "wilkinson@quantumstudies.org" "edu"

# --- Body paragraph 1 ---------------------------------------------------
Replace-Text "Quantum entanglement, an enigmatic phenomenon unveiled by the realm of quantum physics, unravels a world where particles are inextricably linked, challenging our conventional understanding of space, time, and separability" "In the realm of academia, mathematics emerges as a subject of profound significance"
Replace-Text " Delving into the intricate web of entanglement, we discover a profound interconnectedness that transcends physical boundaries" " often hailed as the language of the universe"
Replace-Text " This profound phenomenon, as elucidated by Albert Einstein, is at the heart of a revolution in scientific understanding and exploration" " It unveils the intricate patterns that govern our physical world, epitomizes the essence of critical thinking, and fosters problem solving skills vital to a plethora of disciplines"
Replace-Text " The rich and complex tapestry of entanglement offers a rare glimpse into the fundamental workings of the universe, provoking profound implications for physics, technology, and our understanding of reality itself" " Mathematics is an adventure of exploration, unlocking the enigma of numbers and unraveling the mysteries of shapes. Through its complexities, we find beauty and capture glimpses of the underlying elegant symphony of the universe"

Replace-Text "As we meticulously unravel the threads of entanglement, we unravel secrets long hidden within the fabric of existence" "Mathematics has always intrigued me"
Replace-Text " We discover the remarkable ability of particles to communicate instantaneously, regardless of the vast distances separating them, defying the constraints of space and time" " I still remember the sense of wonder I experienced as a child, playing with blocks and creating intricate patterns"
Replace-Text " This phenomenon, known as quantum teleportation, holds the promise of revolutionizing communication and information processing, laying the groundwork for a new era of technological innovation" " This curiosity has only intensified as I have grown older, and I am fascinated by the way that mathematics can be used to model and understand the world around me"
Replace-Text " The intricate dance of entanglement has illuminated fascinating paradoxes that challenge our classical notions of locality and determinism, inviting us to rethink the very foundation of physics" " I am excited to share my passion for mathematics with my students, and I believe that this subject is essential for helping them to become informed and capable citizens in the 21st century"

Replace-Text "The profound implications of entanglement extend beyond the realm of theoretical physics, reaching into the practical world with remarkable applications" "As a master of numbers, the mathematician wields a powerful tool"
Replace-Text " For instance, the precise control and manipulation of entangled particles hold the key to developing groundbreaking technologies such as quantum computing, encrypted communication, and ultra-sensitive sensors, promising to transform industries and revolutionize sectors" " Numbers, those abstract yet ubiquitous symbols, form the foundation of our universe"
Replace-Text " Moreover, the exploration of entanglement has ignited a profound philosophical discourse on interconnectedness, questioning the nature of reality, consciousness, and our place within the vast cosmos" " They permeate the cosmos, from the arrangement of galaxies to the intricacies of atomic structure"
Replace-Text " It invites us to delve into the profound implications of a universe where everything is fundamentally interconnected, where the actions of one particle can instantaneously influence the behavior of another, regardless of the distance separating them" " Mathematics offers us insights into this cosmos by unlocking the patterns and relationships that bind these numbers, elucidating the symphony that orchestrates the complexity of our physical world"

# --- Summary paragraph ---------------------------------------------------
Replace-Text "Quantum entanglement stands as a testament to the profound interconnectedness that permeates the fabric of reality" "Mathematics is not merely a collection of abstract theories but rather a versatile tool, essential for delving into fields as diverse as engineering and medicine, economics and astronomy"
Replace-Text " Its rich complexity challenges our classical understanding of space, time, and separability, revealing a world where particles dance in harmonious unity" " It cultivates critical thinking, boosts logical reasoning, and hones problem solving abilities"
Replace-Text " Entanglement holds the promise of revolutionizing communication, technology, and our very understanding of the universe. It is a testament to the boundless mysteries that await us, inviting us on an intellectual journey into the heart of reality" " Mathematics is a gateway to understanding the underlying structure of our world, providing budding minds with a formidable weapon in their quest for knowledge"

# --- Trailing empty paragraph added at the end of the body ---------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
